# Rubric-Sprint 3v1.xlsx edit: fill in previously-blank rubric scores for
# the "Other" section (rows 10 and 12) to match the "Possible" column,
# enable iterative calculation, and move the active selection/scroll
# position down to where the user was working (C24 / row 18 visible).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score entry: "Batch process to send emails to rentals overdue" (row 10)
#     and "Report of overdue rentals: 30, 60, 90+ days" (row 12) now earned
#     full points, matching column D ("Possible").
$ws.Range("C10").Value = 5
$ws.Range("C12").Value = 7

# --- Turn on iterative calculation (workbook calcPr -> iterateDelta 1E-4)
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# --- Move the visible/selected cell from C23 down to C24, scrolling so
#     row 18 is at the top of the view.
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select() | Out-Null

$wb.Save() | Out-Null
